$d = $word.ActiveDocument

function Split-RunsAtOffsets($paragraphIndex, $offsets) {
    $p = $d.Paragraphs($paragraphIndex)
    $start = $p.Range.Start
    $names = @()
    $i = 0
    foreach ($off in $offsets) {
        $i = $i + 1
        $markName = "tmpSplitMark_" + $paragraphIndex + "_" + $i
        $r = $d.Range($start + $off, $start + $off)
        $d.Bookmarks.Add($markName, $r)
        $names += $markName
    }
    foreach ($n in $names) {
        $d.Bookmarks($n).Delete()
    }
}

# Paragraph 1 (Title): "Answers: Arithmetic on complex numbers"
# -> "Answers:" " " "Arithmetic" " " "on" " " "complex" " " "numbers"
Split-RunsAtOffsets 1 @(8, 9, 19, 20, 22, 23, 30, 31)

# Paragraph 2 (Author): "Charlotte McCarthy"
# -> "Charlotte" " " "McCarthy"
Split-RunsAtOffsets 2 @(9, 10)

# Paragraph 4 (Abstract): "Answers to questions relating to the guide on arithmetic on complex numbers."
# -> "Answers" " " "to" " " "questions" " " "relating" " " "to" " " "the" " " "guide" " " "on" " " "arithmetic" " " "on" " " "complex" " " "numbers."
Split-RunsAtOffsets 4 @(7, 8, 10, 11, 20, 21, 29, 30, 32, 33, 36, 37, 42, 43, 45, 46, 56, 57, 59, 60, 67, 68)
